$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Cumulative Infections (D) and Current Asymptomatic Infections (E)
# for rows 9-13 to reflect pushing back the start date to Dec 3,
# based on updated infection statistics.

$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 1

$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 4

$ws.Range("D11").Value = 7
$ws.Range("E11").Value = 5

$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 7

$ws.Range("D13").Value = 11
$ws.Range("E13").Value = 9
